$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value = "67ceb73847361997ddc3c691"
$ws.Range("V3").Value = "67ceb73847361997ddc3c68f"
$ws.Range("V4").Value = "67ceb73847361997ddc3c693"
$ws.Range("V5").Value = "67ceb73547361997ddc3c4dd"
$ws.Range("V6").Value = "67ceb73547361997ddc3c4e1"
$ws.Range("V7").Value = "67ceb73547361997ddc3c4df"
$ws.Range("V8").Value = "67ceb73547361997ddc3c4e3"
$ws.Range("V9").Value = "67ceb73947361997ddc3c7a7"
$ws.Range("V10").Value = "67ceb73947361997ddc3c7a9"
$ws.Range("V11").Value = "67ceb73947361997ddc3c7ab"
$ws.Range("V12").Value = "67ceb73c47361997ddc3c8f9"
$ws.Range("V13").Value = "67ceb73c47361997ddc3c8f7"
$ws.Range("V14").Value = "67ceb73747361997ddc3c662"
$ws.Range("V15").Value = "67ceb73747361997ddc3c666"
$ws.Range("V16").Value = "67ceb73747361997ddc3c660"
$ws.Range("V17").Value = "67ceb73747361997ddc3c664"
$ws.Range("V18").Value = "67ceb73a47361997ddc3c801"
$ws.Range("V19").Value = "67ceb73a47361997ddc3c803"
$ws.Range("V20").Value = "67ceb73647361997ddc3c55f"
$ws.Range("V21").Value = "67ceb73647361997ddc3c55d"
$ws.Range("V22").Value = "67ceb73647361997ddc3c55b"
$ws.Range("V23").Value = "67ceb73347361997ddc3c48d"
$ws.Range("V24").Value = "67ceb73347361997ddc3c48f"
$ws.Range("V25").Value = "67ceb73347361997ddc3c491"
$ws.Range("V26").Value = "67ceb73747361997ddc3c633"
$ws.Range("V27").Value = "67ceb73747361997ddc3c631"
$ws.Range("V28").Value = "67ceb73747361997ddc3c637"
$ws.Range("V29").Value = "67ceb73847361997ddc3c6ed"
$ws.Range("V30").Value = "67ceb73847361997ddc3c6f1"
$ws.Range("V31").Value = "67ceb73847361997ddc3c6f7"
$ws.Range("V32").Value = "67ceb73847361997ddc3c6ef"
$ws.Range("V33").Value = "67ceb73847361997ddc3c6f3"
$ws.Range("V34").Value = "67ceb73547361997ddc3c51d"
$ws.Range("V35").Value = "67ceb73547361997ddc3c51f"
$ws.Range("V36").Value = "67ceb73547361997ddc3c51f"
$ws.Range("V37").Value = "67ceb73847361997ddc3c6ba"
$ws.Range("V38").Value = "67ceb73847361997ddc3c6b6"
$ws.Range("V39").Value = "67ceb73847361997ddc3c6b2"
$ws.Range("V40").Value = "67ceb73847361997ddc3c6b4"
$ws.Range("V41").Value = "67ceb73847361997ddc3c6b8"
$ws.Range("V42").Value = "67ceb73a47361997ddc3c7cc"
$ws.Range("V43").Value = "67ceb73a47361997ddc3c7c8"
$ws.Range("V44").Value = "67ceb73a47361997ddc3c7ca"
$ws.Range("V45").Value = "67ceb73a47361997ddc3c7c6"
$ws.Range("V46").Value = "67ceb73947361997ddc3c735"
$ws.Range("V47").Value = "67ceb73947361997ddc3c737"
$ws.Range("V48").Value = "67ceb73947361997ddc3c73b"
$ws.Range("V49").Value = "67ceb73947361997ddc3c739"
$ws.Range("V50").Value = "67ceb73547361997ddc3c508"
$ws.Range("V51").Value = "67ceb73547361997ddc3c50a"
$ws.Range("V52").Value = "67ceb73c47361997ddc3c90c"
$ws.Range("V53").Value = "67ceb73c47361997ddc3c90e"
$ws.Range("V54").Value = "67ceb73c47361997ddc3c910"
$ws.Range("V55").Value = "67ceb73947361997ddc3c77b"
$ws.Range("V56").Value = "67ceb73947361997ddc3c777"
$ws.Range("V57").Value = "67ceb73947361997ddc3c775"
$ws.Range("V58").Value = "67ceb73a47361997ddc3c857"
$ws.Range("V59").Value = "67ceb73a47361997ddc3c855"
$ws.Range("V60").Value = "67ceb73a47361997ddc3c859"
$ws.Range("V61").Value = "67ceb73a47361997ddc3c85b"
$ws.Range("V62").Value = "67ceb73647361997ddc3c5a9"
$ws.Range("V63").Value = "67ceb73647361997ddc3c5ab"
$ws.Range("V64").Value = "67ceb73647361997ddc3c5a7"
$ws.Range("V65").Value = "67ceb73747361997ddc3c5e3"
$ws.Range("V66").Value = "67ceb73747361997ddc3c5e7"
$ws.Range("V67").Value = "67ceb73647361997ddc3c5cc"
$ws.Range("V68").Value = "67ceb73647361997ddc3c5ca"
$ws.Range("V69").Value = "67ceb73b47361997ddc3c890"
$ws.Range("V70").Value = "67ceb73b47361997ddc3c88c"
$ws.Range("V71").Value = "67ceb73b47361997ddc3c888"
$ws.Range("V72").Value = "67ceb73747361997ddc3c604"
$ws.Range("V73").Value = "67ceb73747361997ddc3c600"
$ws.Range("V74").Value = "67ceb73747361997ddc3c606"
$ws.Range("V75").Value = "67ceb73747361997ddc3c602"
$ws.Range("V76").Value = "67ceb73d47361997ddc3c980"
$ws.Range("V77").Value = "67ceb73d47361997ddc3c984"
$ws.Range("V78").Value = "67ceb73d47361997ddc3c982"
$ws.Range("V79").Value = "67ceb73d47361997ddc3c986"
$ws.Range("V80").Value = "67ceb73647361997ddc3c58f"
$ws.Range("V81").Value = "67ceb73d47361997ddc3c9ad"
$ws.Range("V82").Value = "67ceb73d47361997ddc3c9ab"
$ws.Range("V83").Value = "67ceb73d47361997ddc3c9af"
$ws.Range("V84").Value = "67ceb73a47361997ddc3c81a"
$ws.Range("V85").Value = "67ceb73a47361997ddc3c81c"
$ws.Range("V86").Value = "67ceb73a47361997ddc3c833"
$ws.Range("V87").Value = "67ceb73c47361997ddc3c92b"
$ws.Range("V88").Value = "67ceb73a47361997ddc3c84a"
$ws.Range("V89").Value = "67ceb73647361997ddc3c542"
$ws.Range("V90").Value = "67ceb73647361997ddc3c540"
$ws.Range("V91").Value = "67ceb73447361997ddc3c4b4"
$ws.Range("V92").Value = "67ceb73447361997ddc3c4b6"
$ws.Range("V93").Value = "67ceb73447361997ddc3c4b0"
$ws.Range("V94").Value = "67ceb73447361997ddc3c4b2"
$ws.Range("V95").Value = "67ceb73b47361997ddc3c8c8"
$ws.Range("V96").Value = "67ceb73b47361997ddc3c8ca"
$ws.Range("V97").Value = "67ceb73b47361997ddc3c8ce"
$ws.Range("V98").Value = "67ceb73b47361997ddc3c8cc"
$ws.Range("V99").Value = "67ceb73b47361997ddc3c8b5"
$ws.Range("V100").Value = "67ceb73b47361997ddc3c8b3"
$ws.Range("V101").Value = "67ceb73947361997ddc3c796"
$ws.Range("V102").Value = "67ceb73c47361997ddc3c94b"
$ws.Range("V103").Value = "67ceb73c47361997ddc3c953"
$ws.Range("V104").Value = "67ceb73c47361997ddc3c953"
$ws.Range("V105").Value = "67ceb73c47361997ddc3c94d"
$ws.Range("V106").Value = "67ceb73c47361997ddc3c94f"
$ws.Range("V107").Value = "67ceb73c47361997ddc3c93a"
$ws.Range("V108").Value = "67ceb73d47361997ddc3c9d8"
$ws.Range("V109").Value = "67ceb73d47361997ddc3c9d6"
$ws.Range("V110").Value = "67ceb73e47361997ddc3ca1f"
$ws.Range("V111").Value = "67ceb73e47361997ddc3ca23"
$ws.Range("V113").Value = "67ceb73f47361997ddc3cb01"
$ws.Range("V115").Value = "67ceb73e47361997ddc3ca38"
$ws.Range("V116").Value = "67ceb73f47361997ddc3ca98"
$ws.Range("V121").Value = "67ceb73f47361997ddc3cab5"
$ws.Range("V122").Value = "67ceb73f47361997ddc3caa3"
$ws.Range("V127").Value = "67ceb73f47361997ddc3cac5"
$ws.Range("V128").Value = "67ceb73e47361997ddc3ca5a"
$ws.Range("V129").Value = "67ceb73f47361997ddc3cad7"
$ws.Range("V130").Value = "67ceb73f47361997ddc3cae4"
$ws.Range("V131").Value = "67ceb73e47361997ddc3ca0a"
$ws.Range("V132").Value = "67ceb73f47361997ddc3caf4"
$ws.Range("V133").Value = "67ceb73f47361997ddc3ca80"